$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the segment-coordinate parameters for rows 6 and 7 (E6 <-> E7)
$ws.Range("E6").Value = "{'x1':-0.5, 'y1':0, 'x2':0, 'y2':-1}"
$ws.Range("E7").Value = "{'x1':0.5, 'y1':1, 'x2':0, 'y2':-1}"

# Renumber the "paso" column for rows 9-14 (each decreased by 1)
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12

# Fix typo "isoceles" -> "isosceles" in B11
$ws.Range("B11").Value = "Por la afirmación (9) el triángulo `$ \triangle ABC`$ es isósceles y, por tanto, `$\angle ABC \cong \angle BCA`$(Contradicción)"

# Remove parentheses around the angle reference in B13
$ws.Range("B13").Value = "Si `$\overline{AB}`$ es el lado mayor, deducimos que su ángulo opuesto `$ \angle  BCA`$ es mayor y, por tanto, `$ \angle ABC < \angle BCA`$ (Contradicción)"

# Update the selected cell shown in the sheet view
$ws.Range("A15").Select()
